$wb = $excel.ActiveWorkbook

# ---------- Sheet1 ("global") ----------
$ws1 = $wb.Worksheets.Item("global")
$ws1.Range("B4").Value = "localhost"
$ws1.Range("B5").Value = 3306

# ---------- Sheet2 ("dependencies") ----------
$ws2 = $wb.Worksheets.Item("dependencies")
$ws2.Range("K1").Value = "far_responses_bool"
$ws2.Range("L1").Value = "far_responses_date"
$ws2.Range("M1").Value = "far_responses_numeric"
$ws2.Range("N1").Value = "far_responses_options"
$ws2.Range("O1").Value = "far_responses_text"

$ws2.Range("K9").Value = "event-plot"
$ws2.Range("L9").Value = "event-plot"
$ws2.Range("M9").Value = "event-plot"
$ws2.Range("N9").Value = "event-plot"
$ws2.Range("O9").Value = "event-plot"

$ws2.Range("J7").Value = "technical-person"
$ws2.Range("J6").Value = "technical-document"

$ws2.Range("A11").Value = "far_responses_bool"
$ws2.Range("A12").Value = "far_responses_date"
$ws2.Range("A13").Value = "far_responses_numeric"
$ws2.Range("A14").Value = "far_responses_options"
$ws2.Range("A15").Value = "far_responses_text"

$ws2.Range("B29").Select() | Out-Null

# ---------- Sheet3 ("additional") ----------
$ws3 = $wb.Worksheets.Item("additional")
$ws3.Range("A11").Value = "far_responses_bool"
$ws3.Range("B11").Value = 0
$ws3.Range("C11").Value = 0

$ws3.Range("A12").Value = "far_responses_date"
$ws3.Range("B12").Value = 0
$ws3.Range("C12").Value = 0

$ws3.Range("A13").Value = "far_responses_numeric"
$ws3.Range("B13").Value = 0
$ws3.Range("C13").Value = 0

$ws3.Range("A14").Value = "far_responses_options"
$ws3.Range("B14").Value = 0
$ws3.Range("C14").Value = 0

$ws3.Range("A15").Value = "far_responses_text"
$ws3.Range("B15").Value = 0
$ws3.Range("C15").Value = 0

$ws3.Columns.Item(1).ColumnWidth = 20.498697916666668

$ws3.Range("A11:A15").Select() | Out-Null

# ---------- Restore "global" as the active sheet with its final selection ----------
$ws1.Activate() | Out-Null
$ws1.Range("B2:B3").Select() | Out-Null
